$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '42.620.76'
$ws.Range('E2').Value = '  -1.67%  '
$ws.Range('D3').Value = '2.237.93'
$ws.Range('E3').Value = '  -1.86%  '
$ws.Range('E4').Value = '  +0.23%  '
Set-TextValue 'D5' '114.84'
$ws.Range('E5').Value = '  +2.20%  '
Set-TextValue 'D6' '286.08'
$ws.Range('E6').Value = '  +7.77%  '
Set-TextValue 'D7' '0.627'
$ws.Range('E7').Value = '  -3.74%  '
$ws.Range('E8').Value = '  +0.05%  '
Set-TextValue 'D9' '0.613'
$ws.Range('E9').Value = '  +0.29%  '
Set-TextValue 'D10' '46.70'
$ws.Range('E10').Value = '  -0.32%  '
$ws.Range('E11').Value = '  -0.62%  '
Set-TextValue 'D12' '9.13'
$ws.Range('E12').Value = '  -1.59%  '
$ws.Range('E13').Value = '  -3.15%  '
Set-TextValue 'D14' '15.40'
$ws.Range('E14').Value = '  +0.61%  '
Set-TextValue 'D15' '0.882'
$ws.Range('E15').Value = '  +2.08%  '
$ws.Range('D16').Value = '2.576.84'
$ws.Range('D17').Value = '2.236.80'
$ws.Range('E17').Value = '  -1.68%  '
$ws.Range('D18').Value = '42.695.34'
$ws.Range('E18').Value = '  -1.14%  '
$ws.Range('E19').Value = '  -1.05%  '
Set-TextValue 'D20' '6.91'
$ws.Range('E20').Value = '  +2.46%  '
Set-TextValue 'D21' '73.35'
$ws.Range('E21').Value = '  +1.64%  '
Set-TextValue 'D22' '3.17'
$ws.Range('E22').Value = '  +9.59%  '
$ws.Range('E23').Value = '  -2.68%  '
Set-TextValue 'D24' '232.12'
$ws.Range('E24').Value = '  -1.13%  '
Set-TextValue 'D25' '9.20'
$ws.Range('E25').Value = '  -1.87%  '
Set-TextValue 'D26' '12.14'
$ws.Range('E26').Value = '  +6.24%  '
$ws.Range('E27').Value = '  -1.53%  '
Set-TextValue 'D28' '3.95'
$ws.Range('E28').Value = '  -0.73%  '
Set-TextValue 'D29' '40.08'
$ws.Range('E29').Value = '  -2.39%  '
$ws.Range('B30').Value = 'WEMIXToken'
$ws.Range('C30').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D30' '3.30'
$ws.Range('E30').Value = '  -1.14%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D31' '2.23'
$ws.Range('E31').Value = '  -0.25%  '
Set-TextValue 'D32' '175.43'
$ws.Range('E32').Value = '  +1.29%  '
$ws.Range('E33').Value = '  -2.17%  '
Set-TextValue 'D34' '0.0903'
$ws.Range('E34').Value = '  +0.76%  '
Set-TextValue 'D35' '4.62'
$ws.Range('E35').Value = '  +18.97%  '
Set-TextValue 'D36' '5.58'
$ws.Range('E36').Value = '  -0.96%  '
$ws.Range('E37').Value = '  -3.06%  '
Set-TextValue 'D38' '0.0372'
$ws.Range('E38').Value = '  -1.55%  '
Set-TextValue 'D39' '4.62'
$ws.Range('E39').Value = '  -1.41%  '
$ws.Range('E40').Value = '  +1.42%  '
Set-TextValue 'D41' '2.63'
$ws.Range('E41').Value = '  +1.55%  '
Set-TextValue 'D42' '72.78'
$ws.Range('E42').Value = '  -2.04%  '
$ws.Range('E43').Value = '  -5.84%  '
$ws.Range('E44').Value = '  -1.39%  '
$ws.Range('E45').Value = '  +0.55%  '
Set-TextValue 'D46' '1.33'
$ws.Range('E46').Value = '  -2.27%  '
$ws.Range('E47').Value = '  -7.89%  '
$ws.Range('E48').Value = '  +2.62%  '
Set-TextValue 'D49' '8.56'
$ws.Range('E49').Value = '  +0.00%  '
Set-TextValue 'D50' '0.652'
$ws.Range('E50').Value = '  +7.54%  '
Set-TextValue 'D51' '0.474'
$ws.Range('E51').Value = '  +8.72%  '
